# add scene, delete assets, fix error, add sound
$wb = $excel.ActiveWorkbook

# --- add sound: append new respawn points (scene data) to the "Respawn" sheet ---
$ws = $wb.Worksheets.Item("Respawn")

$newRows = @(
    @(5, 6, 75, 4, -60),
    @(6, 6, 40, 4, 55),
    @(7, 7, 0,  4, -20),
    @(8, 8, 0,  4, -36),
    @(9, 9, 0,  4, -23)
)

$startRow = 6
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
}

# --- add scene: move the active tab / selection over to the Respawn sheet ---
$ws.Activate()
$ws.Range("F10").Select()
